$d = $word.ActiveDocument

# 1. Change the "Group2" column header to "Group0"
$t = $d.Tables.Item(1)
$t.Cell(1, 3).Range.Text = "Group0"

# 2. Remove the "Ethnicity - n (%)" section (header row + 5 data rows:
#    Other, Mixed, Asian or Asian British, Black or Black British,
#    White or White British) from the summary table. The parallel
#    "Ethnicity - n" section (without percentages) further down the
#    table is left untouched.
for ($i = 13; $i -ge 8; $i--) {
    $t.Rows.Item($i).Delete()
}
